$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Pspn"
$ws.Cells.Item(2, 3).Value = "Ret"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.042572666666667
$ws.Cells.Item(2, 8).Value = 9.127718
$ws.Cells.Item(2, 9).Value = 0.9827375016055572
$ws.Cells.Item(2, 10).Value = 0.9827375016055572
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.632885
$ws.Cells.Item(2, 14).Value = 10.898655
$ws.Cells.Item(2, 15).Value = 0.4372849566404539
$ws.Cells.Item(2, 16).Value = 0.4372849566404539
$ws.Cells.Item(2, 17).Value = 11.05331660214333
$ws.Cells.Item(2, 18).Value = 99.47984941928999
$ws.Cells.Item(2, 19).Value = 0.4297363257785341
$ws.Cells.Item(2, 20).Value = 0.429736325778534
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Pspn"
$ws.Cells.Item(3, 3).Value = "Ret"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.042572666666667
$ws.Cells.Item(3, 8).Value = 9.127718
$ws.Cells.Item(3, 9).Value = 0.9827375016055572
$ws.Cells.Item(3, 10).Value = 0.9827375016055572
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.884996
$ws.Cells.Item(3, 14).Value = 11.654988
$ws.Cells.Item(3, 15).Value = 0.4676311822169809
$ws.Cells.Item(3, 16).Value = 0.4676311822169809
$ws.Cells.Item(3, 17).Value = 11.82038263970933
$ws.Cells.Item(3, 18).Value = 106.383443757384
$ws.Cells.Item(3, 19).Value = 0.4595586996847689
$ws.Cells.Item(3, 20).Value = 0.4595586996847689
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Pspn"
$ws.Cells.Item(4, 3).Value = "Ret"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.042572666666667
$ws.Cells.Item(4, 8).Value = 9.127718
$ws.Cells.Item(4, 9).Value = 0.9827375016055572
$ws.Cells.Item(4, 10).Value = 0.9827375016055572
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.004340666666666667
$ws.Cells.Item(4, 14).Value = 0.013022
$ws.Cells.Item(4, 15).Value = 0.0005224795816889323
$ws.Cells.Item(4, 16).Value = 0.0005224795816889323
$ws.Cells.Item(4, 17).Value = 0.01320679375511111
$ws.Cells.Item(4, 18).Value = 0.118861143796
$ws.Cells.Item(4, 19).Value = 0.000513460278748898
$ws.Cells.Item(4, 20).Value = 0.000513460278748898
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Pspn"
$ws.Cells.Item(5, 3).Value = "Ret"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.042572666666667
$ws.Cells.Item(5, 8).Value = 9.127718
$ws.Cells.Item(5, 9).Value = 0.9827375016055572
$ws.Cells.Item(5, 10).Value = 0.9827375016055572
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.7855989999999999
$ws.Cells.Item(5, 14).Value = 2.356797
$ws.Cells.Item(5, 15).Value = 0.09456138156087625
$ws.Cells.Item(5, 16).Value = 0.09456138156087625
$ws.Cells.Item(5, 17).Value = 2.390242044360666
$ws.Cells.Item(5, 18).Value = 21.512178399246
$ws.Cells.Item(5, 19).Value = 0.09292901586350533
$ws.Cells.Item(5, 20).Value = 0.09292901586350533
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Pspn"
$ws.Cells.Item(6, 3).Value = "Ret"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.053445
$ws.Cells.Item(6, 8).Value = 0.160335
$ws.Cells.Item(6, 9).Value = 0.01726249839444284
$ws.Cells.Item(6, 10).Value = 0.01726249839444284
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.632885
$ws.Cells.Item(6, 14).Value = 10.898655
$ws.Cells.Item(6, 15).Value = 0.4372849566404539
$ws.Cells.Item(6, 16).Value = 0.4372849566404539
$ws.Cells.Item(6, 17).Value = 0.194159538825
$ws.Cells.Item(6, 18).Value = 1.747435849425
$ws.Cells.Item(6, 19).Value = 0.007548630861919843
$ws.Cells.Item(6, 20).Value = 0.007548630861919842
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Pspn"
$ws.Cells.Item(7, 3).Value = "Ret"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.053445
$ws.Cells.Item(7, 8).Value = 0.160335
$ws.Cells.Item(7, 9).Value = 0.01726249839444284
$ws.Cells.Item(7, 10).Value = 0.01726249839444284
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.884996
$ws.Cells.Item(7, 14).Value = 11.654988
$ws.Cells.Item(7, 15).Value = 0.4676311822169809
$ws.Cells.Item(7, 16).Value = 0.4676311822169809
$ws.Cells.Item(7, 17).Value = 0.20763361122
$ws.Cells.Item(7, 18).Value = 1.86870250098
$ws.Cells.Item(7, 19).Value = 0.00807248253221204
$ws.Cells.Item(7, 20).Value = 0.00807248253221204
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Pspn"
$ws.Cells.Item(8, 3).Value = "Ret"
$ws.Cells.Item(8, 4).Value = "M1"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.053445
$ws.Cells.Item(8, 8).Value = 0.160335
$ws.Cells.Item(8, 9).Value = 0.01726249839444284
$ws.Cells.Item(8, 10).Value = 0.01726249839444284
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.004340666666666667
$ws.Cells.Item(8, 14).Value = 0.013022
$ws.Cells.Item(8, 15).Value = 0.0005224795816889323
$ws.Cells.Item(8, 16).Value = 0.0005224795816889323
$ws.Cells.Item(8, 17).Value = 0.00023198693
$ws.Cells.Item(8, 18).Value = 0.00208788237
$ws.Cells.Item(8, 19).Value = 0.000009019302940034362
$ws.Cells.Item(8, 20).Value = 0.000009019302940034362
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Pspn"
$ws.Cells.Item(9, 3).Value = "Ret"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.053445
$ws.Cells.Item(9, 8).Value = 0.160335
$ws.Cells.Item(9, 9).Value = 0.01726249839444284
$ws.Cells.Item(9, 10).Value = 0.01726249839444284
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.7855989999999999
$ws.Cells.Item(9, 14).Value = 2.356797
$ws.Cells.Item(9, 15).Value = 0.09456138156087625
$ws.Cells.Item(9, 16).Value = 0.09456138156087625
$ws.Cells.Item(9, 17).Value = 0.041986338555
$ws.Cells.Item(9, 18).Value = 0.377877046995
$ws.Cells.Item(9, 19).Value = 0.001632365697370923
$ws.Cells.Item(9, 20).Value = 0.001632365697370923
